# Generate Report for Handoff
# - Update status text from "Handed back: in sync with en-US" to "Ready for handoff"
# - Refresh the handoff timestamps
# - Narrow the "status" columns (they no longer need to fit the long status text)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview "Latest HO Xliff Generate Date" and de-de "Latest Handoff Datetime"
# both shared the same old value (2016-09-03 13:03:58) and move to 13:04:45
$wsOverview.Range("G2").Value = "2016-09-03 13:04:45"
$wsDeDe.Range("H2").Value = "2016-09-03 13:04:45"

# zh-cn "Latest Handoff Datetime" moves from 13:03:54 to 13:04:41
$wsZhCn.Range("H2").Value = "2016-09-03 13:04:41"

# --- Column width changes: the two status-ish columns narrow from ~29.98 to ~17.22 ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
